# PayrollFundTransfer.xlsx update
# Adds a new "Transaction Check" block (T.C (Azure) / T.C (Desc.) / Error) to both
# worksheets, fills in an auto-generated numeric id, and applies a bold 14pt header
# with a boxed (medium outside / thin inside) border around the whole table.

$wb = $excel.ActiveWorkbook

$xlThin = 2
$xlMedium = -4138
$xlLeft = -4131

# ----------------------------------------------------------------------------------
# Sheet 1: PayrollFundTransfer
# ----------------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header labels / value
$ws1.Range("M1").Value = "T.C (Azure)"
$ws1.Range("N1").Value = "T.C (Desc.)"
$ws1.Range("O1").Value = "Error"
$ws1.Range("M2").Value = 119832

# Bold, 14pt header row font
$hdr1 = $ws1.Range("A1:O1")
$hdr1.Font.Size = 14
$hdr1.Font.Bold = $true

# Row heights (visual match for the taller header / bottom row)
$ws1.Rows.Item(1).RowHeight = 18.75
$ws1.Rows.Item(2).RowHeight = 15.75

# ---- Borders -----------------------------------------------------------------
# Interior header cells (B1:N1): thin/thin/medium(top)/thin
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","M","N")) {
  $c = $ws1.Range($col + "1")
  $c.Borders.Weight = $xlThin
  $c.Borders.Item(8).Weight = $xlMedium
}

# Interior bottom-row cells (B2:L2, N2): thin/thin/thin/medium(bottom)
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","N")) {
  $c = $ws1.Range($col + "2")
  $c.Borders.Weight = $xlThin
  $c.Borders.Item(9).Weight = $xlMedium
}

# Top-right corner (O1): thin/medium(right)/medium(top)/thin
$c = $ws1.Range("O1")
$c.Borders.Weight = $xlThin
$c.Borders.Item(8).Weight = $xlMedium
$c.Borders.Item(10).Weight = $xlMedium

# Top-left corner (A1): medium(left)/thin/medium(top)/thin
$c = $ws1.Range("A1")
$c.Borders.Weight = $xlThin
$c.Borders.Item(8).Weight = $xlMedium
$c.Borders.Item(7).Weight = $xlMedium

# Bottom-left corner (A2): medium(left)/thin/thin/medium(bottom)
$c = $ws1.Range("A2")
$c.Borders.Weight = $xlThin
$c.Borders.Item(9).Weight = $xlMedium
$c.Borders.Item(7).Weight = $xlMedium

# Bottom-right corner (O2): thin/medium(right)/thin/medium(bottom)
$c = $ws1.Range("O2")
$c.Borders.Weight = $xlThin
$c.Borders.Item(9).Weight = $xlMedium
$c.Borders.Item(10).Weight = $xlMedium

# M2 (same border as interior bottom row) + left aligned
$c = $ws1.Range("M2")
$c.Borders.Weight = $xlThin
$c.Borders.Item(9).Weight = $xlMedium
$c.HorizontalAlignment = $xlLeft

# ----------------------------------------------------------------------------------
# Sheet 2: PayrollFundTransferAuth
# ----------------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B1").Value = "T.C (Azure)"
$ws2.Range("C1").Value = "T.C (Desc.)"
$ws2.Range("D1").Value = "Error"
$ws2.Range("B2").Value = 119833

$hdr2 = $ws2.Range("A1:D1")
$hdr2.Font.Size = 14
$hdr2.Font.Bold = $true

$ws2.Rows.Item(1).RowHeight = 18.75
$ws2.Rows.Item(2).RowHeight = 15.75

# Interior header cell (B1): thin/thin/medium(top)/thin
$c = $ws2.Range("B1")
$c.Borders.Weight = $xlThin
$c.Borders.Item(8).Weight = $xlMedium

# Interior header cell (C1): same as B1
$c = $ws2.Range("C1")
$c.Borders.Weight = $xlThin
$c.Borders.Item(8).Weight = $xlMedium

# Top-right corner (D1)
$c = $ws2.Range("D1")
$c.Borders.Weight = $xlThin
$c.Borders.Item(8).Weight = $xlMedium
$c.Borders.Item(10).Weight = $xlMedium

# Top-left corner (A1)
$c = $ws2.Range("A1")
$c.Borders.Weight = $xlThin
$c.Borders.Item(8).Weight = $xlMedium
$c.Borders.Item(7).Weight = $xlMedium

# Bottom-left corner (A2)
$c = $ws2.Range("A2")
$c.Borders.Weight = $xlThin
$c.Borders.Item(9).Weight = $xlMedium
$c.Borders.Item(7).Weight = $xlMedium

# Interior bottom-row cell (C2)
$c = $ws2.Range("C2")
$c.Borders.Weight = $xlThin
$c.Borders.Item(9).Weight = $xlMedium

# Bottom-right corner (D2)
$c = $ws2.Range("D2")
$c.Borders.Weight = $xlThin
$c.Borders.Item(9).Weight = $xlMedium
$c.Borders.Item(10).Weight = $xlMedium

# B2 (same border as interior bottom row) + left aligned
$c = $ws2.Range("B2")
$c.Borders.Weight = $xlThin
$c.Borders.Item(9).Weight = $xlMedium
$c.HorizontalAlignment = $xlLeft

$ws2.Range("B2").Select()

# Re-activate sheet 1 (and its own selection) so it remains the visible/selected tab,
# matching the original workbook where PayrollFundTransfer is the tab shown on open.
$ws1.Range("M2").Select()
